$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-period"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# The "Extension.url" row's Fixed Value (Q5) mirrors the same URL text, update it too
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-period"

# The "Extension" row's Constraint(s) cell (AI2) is cleared; that constraint text now only
# belongs to the "Extension.extension" row (AI4), which already carries it.
$elements.Range("AI2").Value = ""
